$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the text "100" in D4 with the number 1000 (rate replaced with new value)
$ws.Range("D4").Value = 1000

# Replace the text "500" in H5 with the number 600 (rate replaced with new value)
$ws.Range("H5").Value = 600

# Update the active selection to D4, reflecting where the user last clicked
$ws.Range("D4").Select()
